$d = $word.ActiveDocument
$wdParagraph = 4
$wdFindContinue = 1

# --- Part 1: remove the "Meta description" paragraph that sits right under the H1 title ---
$metaFind = $d.Content.Duplicate
$found1 = $metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the 'Meta description' paragraph"
}
$metaPara = $d.Range($metaFind.Start, $metaFind.End)
$metaPara.Expand($wdParagraph) | Out-Null
$metaPara.Delete() | Out-Null

# --- Part 2: the closing "Create an eye-catching feature image..." paragraph is replaced by
#             two new paragraphs: a bold re-statement of the page title, followed by the old
#             meta-description text (now rendered in italics) ---
$imgFind = $d.Content.Duplicate
$found2 = $imgFind.Find.Execute("Create an eye-catching feature image", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the 'Create an eye-catching feature image...' paragraph"
}
$imgPara = $d.Range($imgFind.Start, $imgFind.End)
$imgPara.Expand($wdParagraph) | Out-Null
$imgPara.Delete() | Out-Null

$endPoint = $d.Range($d.Content.End, $d.Content.End)

$xmlSnippet = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dream Drop Diamonds Free: Review &amp; Demo | Maximum 10M Jackpot</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Dream Drop Diamonds is a high-variance, 3x3 slot game with Dream Drop Bonus, Diamond Collection and Free Spins. Collect diamonds, and win up to 10,000,000! Play Free Demo.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$endPoint.InsertXML($xmlSnippet) | Out-Null
